$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Add the new test case row (row 11) with shared strings for the new
# username/password pair used by the updated end test case.
$ws.Range("A11").Value = "tony@starkenterprises.com"
$ws.Range("B11").Value = "ironman"
